$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 338. This shifts the existing row 338
# (and everything below it, through the former row 431) down by one,
# copying formatting from the row above into the freshly inserted row.
$ws.Rows(338).Insert()

# Populate the newly inserted row 338 with the new data record.
$ws.Range("A338").Value2 = 3
$ws.Range("B338").Value2 = "Femacal de La Calera"
$ws.Range("C338").Value2 = "Coquimbo"
$ws.Range("D338").Value2 = 44841
$ws.Range("E338").Value2 = 5
$ws.Range("F338").Value2 = 100112043
$ws.Range("G338").Value2 = "Pepino ensalada"
$ws.Range("H338").Value2 = "Sin especificar"
$ws.Range("I338").Value2 = "Primera"
$ws.Range("J338").Value2 = 60
$ws.Range("K338").Value2 = 22000
$ws.Range("L338").Value2 = 22000
$ws.Range("M338").Value2 = 22000
$ws.Range("N338").Value2 = "`$/caja 60 unidades"
$ws.Range("O338").Value2 = "Región de Arica y Parinacota"
$ws.Range("P338").Value2 = 367
$ws.Range("Q338").Value2 = 60
$ws.Range("R338").Value2 = "Hortaliza"
